$d = $word.ActiveDocument

# --- Part 1 -----------------------------------------------------------
# Remove the "Meta description: ..." paragraph that currently sits right
# after the Heading1 title paragraph.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- Part 2 -----------------------------------------------------------
# The last paragraph in the document currently holds the AI image-prompt
# text ("Create a cartoon-style feature image ..."). Replace it with two
# paragraphs:
#   1) a bold paragraph repeating the page title
#   2) an italic paragraph with the meta-description text (the same text
#      that used to live in the paragraph removed above, minus the
#      "Meta description" label)
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

$replacementXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr><w:b/></w:rPr>
              <w:t>Play Ape's Dynasty for Free - Review and Bonus Features</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r/>
            <w:r>
              <w:rPr><w:i/></w:rPr>
              <w:t>Read our review of Ape's Dynasty, an online slot game featuring a free spin bonus and up to 200 spins. Play now for free and unlock various bonuses!</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$lastRange.InsertXML($replacementXml) | Out-Null
